$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New cell values, written in the order that makes the new shared-string
#    table land as: 38 Door Direction Tests, 39 W<, 40 Room Adjacencies,
#    41 Door Adjacencies, 42 Walkway Adjacencies, 43 Test Targets,
#    44 Test Occupied.
# ---------------------------------------------------------------------------
$ws.Range("AG3").Value = "Door Direction Tests"
$ws.Range("V13").Value = "W<"
$ws.Range("AG4").Value = "Room Adjacencies"
$ws.Range("AG5").Value = "Door Adjacencies"
$ws.Range("AG6").Value = "Walkway Adjacencies"
$ws.Range("AG7").Value = "Test Targets"
$ws.Range("AG8").Value = "Test Occupied"

# ---------------------------------------------------------------------------
# 2. Re-use existing (theme based) fills by copying formats from cells that
#    already carry them, so we don't fork new literal-rgb fills.
# ---------------------------------------------------------------------------
# AG4 <- F6 (existing style: solid FFFF2AFF fill)
$ws.Range("F6").Copy()
$ws.Range("AG4").PasteSpecial(-4122)

# AG5 and J17 <- L3 (existing style: theme9 tint -0.25)
$ws.Range("L3").Copy()
$ws.Range("AG5").PasteSpecial(-4122)
$ws.Range("L3").Copy()
$ws.Range("J17").PasteSpecial(-4122)

# AG6 <- AD3 (existing style: theme5)
$ws.Range("AD3").Copy()
$ws.Range("AG6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. New literal fills. Order matters: the first brand-new fill/font
#    combination introduced becomes fill11/style10, etc., so V13 (which pairs
#    the new font with the new yellow fill) must be styled before any plain
#    yellow cell, which must come before the purple cell, which must come
#    before the first cyan cell.
# ---------------------------------------------------------------------------
# V13: Helvetica Neue 10 black text on yellow fill.
$v13 = $ws.Range("V13")
$v13.Font.Name = "Helvetica Neue"
$v13.Font.Size = 10
$v13.Font.Color = 0
$v13.Interior.Color = 65535

# Plain yellow fill cells.
$ws.Range("AG3").Interior.Color = 65535
$ws.Range("T7").Interior.Color = 65535
$ws.Range("C12").Interior.Color = 65535
$ws.Range("L12").Interior.Color = 65535
$ws.Range("R16").Interior.Color = 65535
$ws.Range("AA18").Interior.Color = 65535

# Purple fill cell.
$ws.Range("AG8").Interior.Color = 16748481

# Cyan fill cells.
$ws.Range("Q2").Interior.Color = 15859456
$ws.Range("O9").Interior.Color = 15859456
$ws.Range("S12").Interior.Color = 15859456
$ws.Range("F14").Interior.Color = 15859456
$ws.Range("W14").Interior.Color = 15859456
$ws.Range("AG7").Interior.Color = 15859456

# ---------------------------------------------------------------------------
# 4. Sheet-level bookkeeping to match the edited workbook: new column width,
#    selection and dimension follow automatically from the new cells/col.
# ---------------------------------------------------------------------------
$ws.Columns("AG").ColumnWidth = 36.83203125

$ws.Range("Q2").Select()
